# Apply resume edits: remove all comments (and their range markers), and
# fix a series of wording/typo issues in the body text.

$d = $word.ActiveDocument

# --- 1. Remove every comment (this also removes the commentRangeStart/End
#        and commentReference markers tied to them in the body text). ---
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

# --- 2. Simple text fixes: these runs are not immediately followed by a
#        sibling run sharing identical formatting, so a plain Find/Replace
#        will not trigger unwanted run-coalescing. ---

$d.Content.Find.Execute(
    "Worked with Sensu to create healthchecks with remediation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Worked with Sensu to create health checks with remediation.",
    2)

$d.Content.Find.Execute(
    "Assisted with day to day operations. of production, development and QA environments.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Assisted with day to day operations of production, development and QA environments.",
    2)

$d.Content.Find.Execute(
    "Assisted in build out of new QA environment. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Assisted in the build out a new QA environment. ",
    2)

$d.Content.Find.Execute(
    "Created an tool to allow associates",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Created a tool to allow associates",
    2)

$d.Content.Find.Execute(
    "By Utilizing Ruby on Rails I was able to start anew and rapidly prototype and deploy the requested feature set within 2 weeks.  In addition",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "By Utilizing Ruby on Rails, I started anew and rapidly prototyped and deployed the requested feature set within 2 weeks. In addition",
    2)

$d.Content.Find.Execute(
    "managing customers networks and servers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "managing customer’s networks and servers",
    2)

# --- 3. Edits where the edited run sits right next to a sibling run that
#        shares the exact same rPr. A plain Find/Replace (or Range.Text /
#        Range.Delete) on these causes the engine to coalesce the two runs
#        into one, which the target XML does not do (the following <w:r>
#        is untouched). Work around it by briefly flipping Bold on the
#        range before editing (breaking the "identical formatting"
#        coalescing trigger) and then flipping it back off afterwards -
#        this keeps the two runs distinct while leaving the saved
#        formatting exactly as it was. ---

function Set-RunTextKeepSeparate($doc, $oldText, $newText) {
    $rng = $doc.Content
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $start = $rng.Start
    $end = $rng.End

    $rngBoldOn = $doc.Range($start, $end)
    $rngBoldOn.Bold = 1

    $rngText = $doc.Range($start, $end)
    $rngText.Text = $newText

    $newEnd = $start + $newText.Length
    $rngBoldOff = $doc.Range($start, $newEnd)
    $rngBoldOff.Bold = 0
}

# Once the comments above are deleted, this run becomes directly adjacent
# to the (identically formatted) " as a Systems Administrator ... Team."
# run that used to be separated from it by the comment markers, so the
# whole run's text is replaced in one shot (keeping it as a single run)
# with Bold-toggle protection against auto-coalescing with its neighbor.
Set-RunTextKeepSeparate $d `
    " Managed a mixed Windows/Linux production environment. Daily tasks included mitigating network attacks, monitoring infrastructure, resolving system configuration, performance, and capacity issues. I later transitioned into an engineering role where our team was able to automate a majority of the daily operational tasks, reducing the workload of the operations team as well as decreasing downtime across the environment. " `
    " Managed a mixed Windows/Linux production environment. Daily tasks included mitigating network attacks, monitoring infrastructure, resolving system configuration, performance, and capacity issues. I later transitioned into an engineering role where our team automated most of the daily operational tasks, reducing the workload of the operations team as well as decreasing downtime across the environment. "

Set-RunTextKeepSeparate $d `
    " for a range of enterprise applications.  Responsibilities included application deployment, incident triage, and troubleshooting Linux and Windows environments.  Promoted to System Administrator II after demonstrating an aptitude for tackling complicated infrastructure problems from availability to disaster recovery." `
    "for a range of enterprise applications.  Responsibilities included application deployment, incident triage, and troubleshooting Linux and Windows environments. Promoted to System Administrator II after demonstrating an aptitude for tackling complicated infrastructure problems from availability to disaster recovery."

Set-RunTextKeepSeparate $d "Leveraged agile methodologies " "Leveraged agile methodologies"

Write-Output "done"
